# Demo Verification Script and Display CF
# Updates the "Date" log column (B2:B5) on each scenario sheet with the
# timestamps recorded from the latest verification run.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "PayNowCC";       Cell = "B2"; Value = "Wed Apr 09 01:09:00 IST 2025" },
    @{ Sheet = "PayNowCC";       Cell = "B3"; Value = "Wed Apr 09 01:10:14 IST 2025" },
    @{ Sheet = "PayNowCC";       Cell = "B4"; Value = "Wed Apr 09 01:11:21 IST 2025" },
    @{ Sheet = "PayNowCC";       Cell = "B5"; Value = "Wed Apr 09 01:12:26 IST 2025" },

    @{ Sheet = "PayNowCCSCF";    Cell = "B2"; Value = "Wed Apr 09 01:13:33 IST 2025" },
    @{ Sheet = "PayNowCCSCF";    Cell = "B3"; Value = "Wed Apr 09 01:14:55 IST 2025" },
    @{ Sheet = "PayNowCCSCF";    Cell = "B4"; Value = "Wed Apr 09 01:16:12 IST 2025" },
    @{ Sheet = "PayNowCCSCF";    Cell = "B5"; Value = "Wed Apr 09 01:17:28 IST 2025" },

    @{ Sheet = "PayNowCCDCF";    Cell = "B2"; Value = "Wed Apr 09 01:18:44 IST 2025" },
    @{ Sheet = "PayNowCCDCF";    Cell = "B3"; Value = "Wed Apr 09 01:20:00 IST 2025" },
    @{ Sheet = "PayNowCCDCF";    Cell = "B4"; Value = "Wed Apr 09 01:21:17 IST 2025" },
    @{ Sheet = "PayNowCCDCF";    Cell = "B5"; Value = "Wed Apr 09 01:22:33 IST 2025" },

    @{ Sheet = "NoModifyAmount"; Cell = "B2"; Value = "Wed Apr 09 01:25:58 IST 2025" },

    @{ Sheet = "OverUnderPay";   Cell = "B2"; Value = "Wed Apr 09 01:34:51 IST 2025" },
    @{ Sheet = "OverUnderPay";   Cell = "B3"; Value = "Wed Apr 09 01:35:49 IST 2025" },

    @{ Sheet = "NoOverPay";      Cell = "B2"; Value = "Wed Apr 09 01:42:36 IST 2025" }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $ws.Range($u.Cell).Value = $u.Value
}
